# Auto-generated edit script: updates Sheets market-data snapshot values
# per the "chore: update Sheets via scheduled runner" commit.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 811.8570999999999
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H19").Value = 366.34784
$ws.Range("I19").Value = 351.58334
$ws.Range("K19").Value = 351.58334
$ws.Range("M19").Value = -176.58334
$ws.Range("H40").Value = 5018.2
$ws.Range("I40").Value = 4074.75
$ws.Range("J40").Value = 6433.375
$ws.Range("K40").Value = 4074.75
$ws.Range("L40").Value = 6433.375
$ws.Range("M40").Value = -3899.75
$ws.Range("N40").Value = -6783.375
$ws.Range("H129").Value = 2176.4119
$ws.Range("I129").Value = 1759.625
$ws.Range("J129").Value = 2546.889
$ws.Range("K129").Value = 5278.875
$ws.Range("L129").Value = 7640.667
$ws.Range("M129").Value = -278.875
$ws.Range("N129").Value = -17640.667
$ws.Range("H132").Value = 15865.875
$ws.Range("I132").Value = 17181.154
$ws.Range("J132").Value = 10166.333
$ws.Range("K132").Value = 51543.462
$ws.Range("L132").Value = 30498.999
$ws.Range("M132").Value = -49013.462
$ws.Range("N132").Value = -35558.999
$ws.Range("H135").Value = 2883.8572
$ws.Range("I135").Value = 2935.75
$ws.Range("K135").Value = 26421.75
$ws.Range("M135").Value = -23886.75
$ws.Range("H141").Value = 2377.3333
$ws.Range("I141").Value = 1240.75
$ws.Range("K141").Value = 3722.25
$ws.Range("M141").Value = 1457.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5196.486
$ws.Range("I32").Value = 4834.6177
$ws.Range("K32").Value = 4834.6177
$ws.Range("M32").Value = -4547.6177
$ws.Range("H74").Value = 1670
$ws.Range("I74").Value = 1664
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 1664
$ws.Range("L74").Value = 1700
$ws.Range("M74").Value = -790
$ws.Range("N74").Value = -3448
$ws.Range("H77").Value = 1670
$ws.Range("I77").Value = 1664
$ws.Range("J77").Value = 1700
$ws.Range("K77").Value = 8320
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = -3952
$ws.Range("N77").Value = -17236
$ws.Range("H88").Value = 4300
$ws.Range("I88").Value = 7500
$ws.Range("J88").Value = 2700
$ws.Range("K88").Value = 7500
$ws.Range("L88").Value = 2700
$ws.Range("M88").Value = -7094
$ws.Range("N88").Value = -3512
$ws.Range("H91").Value = 4300
$ws.Range("I91").Value = 7500
$ws.Range("J91").Value = 2700
$ws.Range("K91").Value = 7500
$ws.Range("L91").Value = 2700
$ws.Range("M91").Value = -6096
$ws.Range("N91").Value = -5508

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 396.75
$ws.Range("J64").Value = 522.1429000000001
$ws.Range("L64").Value = 522.1429000000001
$ws.Range("N64").Value = -972.1429000000001
$ws.Range("H67").Value = 396.75
$ws.Range("J67").Value = 522.1429000000001
$ws.Range("L67").Value = 522.1429000000001
$ws.Range("N67").Value = -2082.1429
$ws.Range("H86").Value = 7423
$ws.Range("I86").Value = 5996.8887
$ws.Range("K86").Value = 5996.8887
$ws.Range("M86").Value = -4873.8887
$ws.Range("H89").Value = 7423
$ws.Range("I89").Value = 5996.8887
$ws.Range("K89").Value = 29984.4435
$ws.Range("M89").Value = -24368.4435

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 99995
$ws.Range("J9").Value = 99995
$ws.Range("L9").Value = 99995
$ws.Range("N9").Value = -100331
$ws.Range("H16").Value = 2030
$ws.Range("I16").Value = 1887.5
$ws.Range("K16").Value = 1887.5
$ws.Range("M16").Value = -1600.5
$ws.Range("H31").Value = 4069.5305
$ws.Range("I31").Value = 1729.1143
$ws.Range("K31").Value = 1729.1143
$ws.Range("M31").Value = -1434.1143
$ws.Range("H34").Value = 4069.5305
$ws.Range("I34").Value = 1729.1143
$ws.Range("K34").Value = 1729.1143
$ws.Range("M34").Value = -1527.1143
$ws.Range("H35").Value = 225.57143
$ws.Range("I35").Value = 225.57143
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 225.57143
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 68.42857000000001
$ws.Range("N35").ClearContents()
$ws.Range("H39").Value = 2432
$ws.Range("I39").Value = 2432
$ws.Range("K39").Value = 2432
$ws.Range("M39").Value = -2041
$ws.Range("H49").Value = 2432
$ws.Range("I49").Value = 2432
$ws.Range("K49").Value = 2432
$ws.Range("M49").Value = -2250
$ws.Range("H58").Value = 4394.5356
$ws.Range("I58").Value = 2448.1304
$ws.Range("K58").Value = 2448.1304
$ws.Range("M58").Value = -2245.1304
$ws.Range("H99").Value = 3413.4482
$ws.Range("I99").Value = 3375.762
$ws.Range("K99").Value = 3375.762
$ws.Range("M99").Value = -1877.762
$ws.Range("H113").Value = 2030
$ws.Range("I113").Value = 1887.5
$ws.Range("K113").Value = 1887.5
$ws.Range("M113").Value = 282.5
$ws.Range("H122").Value = 1755.1875
$ws.Range("I122").Value = 1714.1428
$ws.Range("J122").Value = 2042.5
$ws.Range("K122").Value = 5142.428400000001
$ws.Range("L122").Value = 6127.5
$ws.Range("M122").Value = -2692.428400000001
$ws.Range("N122").Value = -11027.5
$ws.Range("H126").Value = 3413.4482
$ws.Range("I126").Value = 3375.762
$ws.Range("K126").Value = 10127.286
$ws.Range("M126").Value = -7657.286
$ws.Range("H132").Value = 5311.125
$ws.Range("I132").Value = 4642
$ws.Range("K132").Value = 13926
$ws.Range("M132").Value = -11396
$ws.Range("H136").Value = 4394.5356
$ws.Range("I136").Value = 2448.1304
$ws.Range("K136").Value = 7344.3912
$ws.Range("M136").Value = -4794.3912

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 7115.3335
$ws.Range("I139").Value = 6531.8
$ws.Range("K139").Value = 19595.4
$ws.Range("M139").Value = -14455.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3065.889
$ws.Range("I80").Value = 3019.4
$ws.Range("J80").Value = 3124
$ws.Range("K80").Value = 3019.4
$ws.Range("L80").Value = 3124
$ws.Range("M80").Value = -2021.4
$ws.Range("N80").Value = -5120
$ws.Range("H83").Value = 3065.889
$ws.Range("I83").Value = 3019.4
$ws.Range("J83").Value = 3124
$ws.Range("K83").Value = 15097
$ws.Range("L83").Value = 15620
$ws.Range("M83").Value = -10105
$ws.Range("N83").Value = -25604
$ws.Range("H97").Value = 938.53845
$ws.Range("I97").Value = 882.36365
$ws.Range("J97").Value = 1247.5
$ws.Range("K97").Value = 882.36365
$ws.Range("L97").Value = 1247.5
$ws.Range("M97").Value = -386.36365
$ws.Range("N97").Value = -2239.5
$ws.Range("H102").Value = 2213.0417
$ws.Range("I102").Value = 1814.9524
$ws.Range("K102").Value = 1814.9524
$ws.Range("M102").Value = -192.9523999999999
$ws.Range("H113").Value = 9362.5
$ws.Range("I113").Value = 7450
$ws.Range("K113").Value = 7450
$ws.Range("M113").Value = -5280
$ws.Range("H132").Value = 32446.695
$ws.Range("I132").Value = 36588.324
$ws.Range("J132").Value = 6768.6
$ws.Range("K132").Value = 109764.972
$ws.Range("L132").Value = 20305.8
$ws.Range("M132").Value = -107234.972
$ws.Range("N132").Value = -25365.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 319.8
$ws.Range("I16").Value = 319.8
$ws.Range("K16").Value = 319.8
$ws.Range("M16").Value = -149.8
$ws.Range("H32").Value = 1671
$ws.Range("I32").Value = 1671
$ws.Range("K32").Value = 1671
$ws.Range("M32").Value = -1354
$ws.Range("H40").Value = 6450.6665
$ws.Range("I40").Value = 5565
$ws.Range("J40").Value = 9993.333000000001
$ws.Range("K40").Value = 5565
$ws.Range("L40").Value = 9993.333000000001
$ws.Range("M40").Value = -5429
$ws.Range("N40").Value = -10265.333
$ws.Range("H55").Value = 1257.7142
$ws.Range("I55").Value = 1538.8334
$ws.Range("K55").Value = 1538.8334
$ws.Range("M55").Value = -1365.8334
$ws.Range("H61").Value = 3199.6538
$ws.Range("I61").Value = 2326.8635
$ws.Range("K61").Value = 2326.8635
$ws.Range("M61").Value = -2124.8635
$ws.Range("H100").Value = 5335.524
$ws.Range("I100").Value = 1838.4445
$ws.Range("J100").Value = 7958.3335
$ws.Range("K100").Value = 1838.4445
$ws.Range("L100").Value = 7958.3335
$ws.Range("M100").Value = -1297.4445
$ws.Range("N100").Value = -9040.333500000001
$ws.Range("H113").Value = 3199.6538
$ws.Range("I113").Value = 2326.8635
$ws.Range("K113").Value = 2326.8635
$ws.Range("M113").Value = -156.8634999999999
$ws.Range("H122").Value = 3961.3333
$ws.Range("I122").Value = 3961.3333
$ws.Range("K122").Value = 11883.9999
$ws.Range("M122").Value = -9433.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 10605000
$ws.Range("J94").Value = 10605000
$ws.Range("L94").Value = 10605000
$ws.Range("N94").Value = -10606802
$ws.Range("H122").Value = 2533.45
$ws.Range("I122").Value = 2275.7273
$ws.Range("J122").Value = 2848.4443
$ws.Range("K122").Value = 6827.1819
$ws.Range("L122").Value = 8545.332900000001
$ws.Range("M122").Value = -4377.1819
$ws.Range("N122").Value = -13445.3329
$ws.Range("H132").Value = 3387.8572
$ws.Range("I132").Value = 3233.2
$ws.Range("J132").Value = 3774.5
$ws.Range("K132").Value = 9699.599999999999
$ws.Range("L132").Value = 11323.5
$ws.Range("M132").Value = -7169.599999999999
$ws.Range("N132").Value = -16383.5
$ws.Range("H136").Value = 4455.7085
$ws.Range("I136").Value = 2997.9
$ws.Range("J136").Value = 5497
$ws.Range("K136").Value = 8993.700000000001
$ws.Range("L136").Value = 16491
$ws.Range("M136").Value = -6443.700000000001
$ws.Range("N136").Value = -21591
$ws.Range("H141").Value = 179949.6
$ws.Range("J141").Value = 212437
$ws.Range("L141").Value = 212437
$ws.Range("N141").Value = -222797

